{"js": "// Merge the split \"ffmpeg -stream_loop ...\" runs back into a single run,\n// and insert the missing \"-f rtsp \" format flag before the rtsp:// URL.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst newText =\n  \"ffmpeg -stream_loop -1 -re -i case1.mp4 -c copy -f rtsp rtsp://192.168.1.177:8554/mystream\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"ffmpeg -stream_loop\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the ffmpeg command paragraph.\");\n}\n\n// Replacing the paragraph's text collapses all of its runs into one run\n// while keeping the formatting (rPr) of the paragraph's first run.\ntarget.insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Merge the split \"ffmpeg -stream_loop ...\" runs back into a single run,\n# and insert the missing \"-f rtsp \" format flag before the rtsp:// URL.\n$d = $word.ActiveDocument\n\n$oldText = \"ffmpeg -stream_loop -1 -re -i case1.mp4 -c copy -f rtsp://192.168.1.177:8554/mystream\"\n$newText = \"ffmpeg -stream_loop -1 -re -i case1.mp4 -c copy -f rtsp rtsp://192.168.1.177:8554/mystream\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$replaced = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\nif (-not $replaced) {\n    throw \"Could not find the ffmpeg command text to replace.\"\n}\n"}
